$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "45.932.71"
$ws.Range("E2").Value = "  -1.61%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.396.51"
$ws.Range("E3").Value = "  +3.89%  "

# Row 4
$ws.Range("E4").Value = "  -0.13%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "302.33"
$ws.Range("E5").Value = "  +0.67%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "98.28"
$ws.Range("E6").Value = "  +0.76%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.568"
$ws.Range("E7").Value = "  -0.65%  "

# Row 8
$ws.Range("E8").Value = "  +0.05%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.516"
$ws.Range("E9").Value = "  -3.19%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "34.51"
$ws.Range("E10").Value = "  -3.17%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0798"
$ws.Range("E11").Value = "  -0.80%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.16"
$ws.Range("E12").Value = "  -3.54%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.104"
$ws.Range("E13").Value = "  -0.30%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.749.69"
$ws.Range("E14").Value = "  +3.53%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.369.94"
$ws.Range("E15").Value = "  +2.66%  "

# Row 16
$ws.Range("B16").Value = "Polygon"
$ws.Range("C16").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.824"
$ws.Range("E16").Value = "  -0.69%  "

# Row 17
$ws.Range("B17").Value = "Chainlink"
$ws.Range("C17").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.87"
$ws.Range("E17").Value = "  -0.96%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "45.959.91"
$ws.Range("E18").Value = "  -1.20%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.79"
$ws.Range("E19").Value = "  -3.64%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0977"
$ws.Range("E20").Value = "  +3.29%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.13"
$ws.Range("E21").Value = "  -0.92%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "67.17"
$ws.Range("E22").Value = "  +0.15%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "245.33"
$ws.Range("E23").Value = "  -1.90%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.84"
$ws.Range("E24").Value = "  -3.45%  "

# Row 25
$ws.Range("E25").Value = "  -0.02%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.93"
$ws.Range("E26").Value = "  -3.31%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "38.89"
$ws.Range("E27").Value = "  -7.64%  "

# Row 28
$ws.Range("E28").Value = "  -3.74%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.81"
$ws.Range("E29").Value = "  -0.57%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "21.76"
$ws.Range("E30").Value = "  +8.09%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.78"
$ws.Range("E31").Value = "  +21.86%  "

# Row 32
$ws.Range("E32").Value = "  +7.16%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.59"
$ws.Range("E33").Value = "  -3.84%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "144.55"
$ws.Range("E34").Value = "  -0.71%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0775"
$ws.Range("E35").Value = "  -3.86%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.113"
$ws.Range("E36").Value = "  +1.10%  "

# Row 37
$ws.Range("B37").Value = "ARBITRUM"
$ws.Range("C37").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.93"
$ws.Range("E37").Value = "  +6.70%  "

# Row 38
$ws.Range("B38").Value = "Stellar"
$ws.Range("C38").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.117"
$ws.Range("E38").Value = "  -2.05%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "14.85"
$ws.Range("E39").Value = "  -0.69%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.92"
$ws.Range("E40").Value = "  -3.43%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0301"
$ws.Range("E41").Value = "  -2.49%  "

# Row 42
$ws.Range("E42").Value = "  -6.59%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.999"
$ws.Range("E43").Value = "  +0.04%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.873.25"
$ws.Range("E44").Value = "  +2.57%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "88.96"
$ws.Range("E45").Value = "  -3.07%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.76"
$ws.Range("E46").Value = "  -10.30%  "

# Row 47
$ws.Range("B47").Value = "FraxShare"
$ws.Range("C47").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.47"
$ws.Range("E47").Value = "  +6.38%  "

# Row 48
$ws.Range("B48").Value = "Algorand"
$ws.Range("C48").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.187"
$ws.Range("E48").Value = "  -4.49%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.620.69"
$ws.Range("E49").Value = "  +3.29%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "98.57"
$ws.Range("E50").Value = "  +0.87%  "

# Row 51
$ws.Range("B51").Value = "THORChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.77"
$ws.Range("E51").Value = "  -3.30%  "
